# "Manger Screen and commenting code"
# Rebuild the "Top topics raised" table: the underlying raw data changed
# (some topics/comments such as "word", "OneNote", "Word", "publisher",
# the standalone "1", the Publisher "not working" comments and
# "PowerPoint not working" dropped out, while the remaining topics'
# rank/order shifted) so the PROBLEM / Number Of Times Raised table
# shrinks from 29 data rows to 21 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the rows that are no longer part of the (shorter) report.
$ws.Rows("23:30").Delete()

# New ordered (PROBLEM, Number Of Times Raised) data for rows 2..22.
$topics = @(
  @("Microsoft Sway ended with an error is not able to open charts", 2),
  @("Microsoft OneNote ended with an error is taking too long to open", 2),
  @("Microsoft Power BI ended with an error is not able to start", 1),
  @("Microsoft OneNote ended with an error not able to start", 1),
  @("Microsoft Office ended with an error is hanging", 1),
  @("Microsoft Office ended with an error not able to start", 1),
  @("Ishan is not able to run Microsoft Excel", 1),
  @("Microsoft Skype ended with an error is not able to connect to API", 1),
  @("Microsoft OneDrive ended with an error is taking too long to open", 1),
  @("Microsoft PictureMgr ended with an error not able to start", 1),
  @("Ishan is now facing issues with Access", 1),
  @("Microsoft Excel ended with an error is not able to open charts", 1),
  @("Microsoft Office ended with an error is not able to start", 1),
  @("Microsoft Excel ended with an error is hanging", 1),
  @("Microsoft Publisher ended with an error is not able to start", 1),
  @("Hi I am not able to open Excel", 1),
  @("Microsoft OneNote ended with an error is not able to connect to API", 1),
  @("Microsoft OneNote ended with an error is not able to start", 1),
  @("Microsoft Project ended with an error is not able to connect to API", 1),
  @("Microsoft Access ended with an error is not able to open charts", 1),
  @("Microsoft PictureMgr ended with an error is hanging", 1)
)

$r = 2
foreach ($topic in $topics) {
  $ws.Cells.Item($r, 1).Value = $topic[0]
  $ws.Cells.Item($r, 2).Value = $topic[1]
  $r = $r + 1
}
